$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6444444444444445
$ws.Range("D2").Value = 0.6987951807228916
$ws.Range("B3").Value = 0.8125
$ws.Range("C3").Value = 0.7090909090909091
$ws.Range("D3").Value = 0.7572815533980584
$ws.Range("B4").Value = 0.7311827956989247
$ws.Range("C4").Value = 0.7311827956989247
$ws.Range("D4").Value = 0.7311827956989247
$ws.Range("E4").Value = 0.7311827956989247
$ws.Range("B5").Value = 0.7284722222222222
$ws.Range("C5").Value = 0.7361244019138756
$ws.Range("D5").Value = 0.728038367060475
$ws.Range("B6").Value = 0.7438321385902031
$ws.Range("C6").Value = 0.7311827956989247
$ws.Range("D6").Value = 0.7333838957458397
$ws.Range("B7").Value = 0.6590909090909091
$ws.Range("C7").Value = 0.7631578947368421
$ws.Range("D7").Value = 0.7073170731707317
$ws.Range("B8").Value = 0.8163265306122449
$ws.Range("C8").Value = 0.7272727272727273
$ws.Range("D8").Value = 0.7692307692307693
$ws.Range("B9").Value = 0.7419354838709677
$ws.Range("C9").Value = 0.7419354838709677
$ws.Range("D9").Value = 0.7419354838709677
$ws.Range("E9").Value = 0.7419354838709677
$ws.Range("B10").Value = 0.737708719851577
$ws.Range("C10").Value = 0.7452153110047848
$ws.Range("D10").Value = 0.7382739212007505
$ws.Range("B11").Value = 0.7520797175175055
$ws.Range("C11").Value = 0.7419354838709677
$ws.Range("D11").Value = 0.7439326998729044
$ws.Range("B12").Value = 0.6046511627906976
$ws.Range("C12").Value = 0.6842105263157895
$ws.Range("D12").Value = 0.6419753086419753
$ws.Range("B13").Value = 0.76
$ws.Range("C13").Value = 0.6909090909090909
$ws.Range("D13").Value = 0.7238095238095238
$ws.Range("B14").Value = 0.6881720430107527
$ws.Range("C14").Value = 0.6881720430107527
$ws.Range("D14").Value = 0.6881720430107527
$ws.Range("E14").Value = 0.6881720430107527
$ws.Range("B15").Value = 0.6823255813953488
$ws.Range("C15").Value = 0.6875598086124401
$ws.Range("D15").Value = 0.6828924162257495
$ws.Range("B16").Value = 0.696524131032758
$ws.Range("C16").Value = 0.6881720430107527
$ws.Range("D16").Value = 0.6903718875045041
$ws.Range("B17").Value = 0.5952380952380952
$ws.Range("D17").Value = 0.625
$ws.Range("B18").Value = 0.7450980392156863
$ws.Range("C18").Value = 0.6909090909090909
$ws.Range("D18").Value = 0.7169811320754716
$ws.Range("B19").Value = 0.6774193548387096
$ws.Range("C19").Value = 0.6774193548387096
$ws.Range("D19").Value = 0.6774193548387096
$ws.Range("E19").Value = 0.6774193548387096
$ws.Range("B20").Value = 0.6701680672268908
$ws.Range("C20").Value = 0.6744019138755981
$ws.Range("D20").Value = 0.6709905660377358
$ws.Range("B21").Value = 0.6838649438269933
$ws.Range("C21").Value = 0.6774193548387096
$ws.Range("D21").Value = 0.6793974437005478
$ws.Range("B22").Value = 0.6153846153846154
$ws.Range("C22").Value = 0.631578947368421
$ws.Range("D22").Value = 0.6233766233766234
$ws.Range("B23").Value = 0.7407407407407407
$ws.Range("C23").Value = 0.7272727272727273
$ws.Range("D23").Value = 0.7339449541284404
$ws.Range("B24").Value = 0.6881720430107527
$ws.Range("C24").Value = 0.6881720430107527
$ws.Range("D24").Value = 0.6881720430107527
$ws.Range("E24").Value = 0.6881720430107527
$ws.Range("B25").Value = 0.6780626780626781
$ws.Range("C25").Value = 0.6794258373205742
$ws.Range("D25").Value = 0.6786607887525319
$ws.Range("B26").Value = 0.6895199583371626
$ws.Range("C26").Value = 0.6881720430107527
$ws.Range("D26").Value = 0.6887664964018915
